# Fix IFRS financial figures for 동성코퍼레이션 (company_list sheet)
# The previously-entered values were wrong (likely pulled from the wrong
# company/scale); this corrects rows 2-6 to the right figures and removes the
# bogus forecast rows 7-9 (their data columns were fabricated placeholders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 9054
$ws.Range("E2").Value = 500
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 314
$ws.Range("H2").Value = 291
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 218
$ws.Range("K2").Value = 8092
$ws.Range("L2").Value = 4308
$ws.Range("M2").Value = 3783
$ws.Range("N2").Value = 2029
$ws.Range("O2").Value = 1754
$ws.Range("P2").Value = 362
$ws.Range("Q2").Value = 623
$ws.Range("R2").Value = -460
$ws.Range("S2").Value = -356
$ws.Range("T2").Value = 289
$ws.Range("U2").Value = 335
$ws.Range("V2").Value = 2388
$ws.Range("W2").Value = 5.53
$ws.Range("X2").Value = 3.21
$ws.Range("Y2").Value = 3.66
$ws.Range("Z2").Value = 3.77
$ws.Range("AA2").Value = 113.88
$ws.Range("AB2").Value = 460.45
$ws.Range("AC2").Value = 203
$ws.Range("AD2").Value = 32.9
$ws.Range("AE2").Value = 5604
$ws.Range("AF2").Value = 1.19
$ws.Range("AG2").Value = 170
$ws.Range("AH2").Value = 2.54
$ws.Range("AI2").Value = 84.52
$ws.Range("AJ2").Value = 36214645

# Row 3
$ws.Range("D3").Value = 8712
$ws.Range("E3").Value = 667
$ws.Range("F3").Value = 667
$ws.Range("G3").Value = 573
$ws.Range("H3").Value = 436
$ws.Range("I3").Value = 203
$ws.Range("J3").Value = 233
$ws.Range("K3").Value = 8017
$ws.Range("L3").Value = 4032
$ws.Range("M3").Value = 3985
$ws.Range("N3").Value = 2663
$ws.Range("O3").Value = 1322
$ws.Range("P3").Value = 454
$ws.Range("Q3").Value = 497
$ws.Range("R3").Value = -333
$ws.Range("S3").Value = -202
$ws.Range("T3").Value = 325
$ws.Range("U3").Value = 171
$ws.Range("V3").Value = 2424
$ws.Range("W3").Value = 7.66
$ws.Range("X3").Value = 5
$ws.Range("Y3").Value = 8.65
$ws.Range("Z3").Value = 5.41
$ws.Range("AA3").Value = 101.18
$ws.Range("AB3").Value = 521.5599999999999
$ws.Range("AC3").Value = 496
$ws.Range("AD3").Value = 12.87
$ws.Range("AE3").Value = 6115
$ws.Range("AF3").Value = 1.04
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 3.13
$ws.Range("AI3").Value = 42.93
$ws.Range("AJ3").Value = 45444970

# Row 4
$ws.Range("D4").Value = 8529
$ws.Range("E4").Value = 650
$ws.Range("F4").Value = 650
$ws.Range("G4").Value = 704
$ws.Range("H4").Value = 530
$ws.Range("I4").Value = 375
$ws.Range("J4").Value = 155
$ws.Range("K4").Value = 8721
$ws.Range("L4").Value = 4314
$ws.Range("M4").Value = 4406
$ws.Range("N4").Value = 2953
$ws.Range("O4").Value = 1454
$ws.Range("P4").Value = 454
$ws.Range("Q4").Value = 783
$ws.Range("R4").Value = -497
$ws.Range("S4").Value = -178
$ws.Range("T4").Value = 291
$ws.Range("U4").Value = 492
$ws.Range("V4").Value = 2392
$ws.Range("W4").Value = 7.63
$ws.Range("X4").Value = 6.21
$ws.Range("Y4").Value = 13.35
$ws.Range("Z4").Value = 6.33
$ws.Range("AA4").Value = 97.90000000000001
$ws.Range("AB4").Value = 585.28
$ws.Range("AC4").Value = 824
$ws.Range("AD4").Value = 7.52
$ws.Range("AE4").Value = 6744
$ws.Range("AF4").Value = 0.92
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 3.23
$ws.Range("AI4").Value = 23.37
$ws.Range("AJ4").Value = 45444970

# Row 5
$ws.Range("D5").Value = 8273
$ws.Range("E5").Value = 37
$ws.Range("F5").Value = 37
$ws.Range("G5").Value = 14
$ws.Range("H5").Value = -59
$ws.Range("I5").Value = 98
$ws.Range("J5").Value = -157
$ws.Range("K5").Value = 8547
$ws.Range("L5").Value = 4368
$ws.Range("M5").Value = 4178
$ws.Range("N5").Value = 2956
$ws.Range("O5").Value = 1222
$ws.Range("P5").Value = 454
$ws.Range("Q5").Value = 472
$ws.Range("R5").Value = -363
$ws.Range("S5").Value = -32
$ws.Range("T5").Value = 488
$ws.Range("U5").Value = -16
$ws.Range("V5").Value = 2637
$ws.Range("W5").Value = 0.44
$ws.Range("X5").Value = -0.72
$ws.Range("Y5").Value = 3.32
$ws.Range("Z5").Value = -0.6899999999999999
$ws.Range("AA5").Value = 104.55
$ws.Range("AB5").Value = 580.39
$ws.Range("AC5").Value = 216
$ws.Range("AD5").Value = 25.98
$ws.Range("AE5").Value = 6719
$ws.Range("AF5").Value = 0.83
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 3.57
$ws.Range("AI5").Value = 89.67
$ws.Range("AJ5").Value = 45444970

# Row 6
$ws.Range("D6").Value = 8131
$ws.Range("E6").Value = 78
$ws.Range("F6").Value = 78
$ws.Range("G6").Value = 29
$ws.Range("H6").Value = -48
$ws.Range("I6").Value = -32
$ws.Range("K6").Value = 8081
$ws.Range("L6").Value = 4074
$ws.Range("M6").Value = 4007
$ws.Range("N6").Value = 2828
$ws.Range("P6").Value = 454
$ws.Range("Q6").Value = 477
$ws.Range("R6").Value = -189
$ws.Range("S6").Value = -404
$ws.Range("T6").Value = 177
$ws.Range("U6").Value = 300
$ws.Range("V6").Value = 2347
$ws.Range("W6").Value = 0.96
$ws.Range("X6").Value = -0.59
$ws.Range("Y6").Value = -1.11
$ws.Range("Z6").Value = -0.58
$ws.Range("AA6").Value = 101.69
$ws.Range("AB6").Value = 554.6
$ws.Range("AC6").Value = -71
$ws.Range("AD6").Value = -72.42
$ws.Range("AE6").Value = 6428
$ws.Range("AF6").Value = 0.8
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 3.9
$ws.Range("AI6").Value = -273.32
$ws.Range("AJ6").Value = 45444970

# Rows 7-9 previously carried a full data block (D:AJ) that does not belong
# here; clear it so only the label columns (A/B/C) remain, matching rows
# as originally intended.
$ws.Range("D7:AJ9").ClearContents()

